# This script rearranges (permutes) the data in rows 2-7 of the active sheet.
# Only the columns A, B, E, F, G, H, Q and R carry row-specific values that
# move between rows; all other columns remain identical across these rows.
#
# The row permutation (target row <- source row) is:
#   2 <- 4
#   3 <- 5
#   4 <- 2
#   5 <- 6
#   6 <- 7
#   7 <- 3

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A", "B", "E", "F", "G", "H", "Q", "R")

# Snapshot the current (pre-edit) values for the affected rows/columns
# before writing anything, since several rows read from one another.
$snapshot = @{}
foreach ($row in 2..7) {
    $rowValues = @{}
    foreach ($col in $cols) {
        $rowValues[$col] = $ws.Range("$col$row").Value2
    }
    $snapshot[$row] = $rowValues
}

# target row -> source row
$mapping = @{
    2 = 4
    3 = 5
    4 = 2
    5 = 6
    6 = 7
    7 = 3
}

foreach ($targetRow in 2..7) {
    $sourceRow = $mapping[$targetRow]
    $sourceValues = $snapshot[$sourceRow]
    foreach ($col in $cols) {
        $ws.Range("$col$targetRow").Value = $sourceValues[$col]
    }
}
